$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (student 204897687): score moved from 0.0/100.0 to 33.33/100.0,
# and the review notes now list the failed checks instead of being blank.
$ws.Range("B2").Value = "33.33 / 100.0"
$ws.Range("C2").Value = "display_name_with_only_last_name: failed`ndisplay_full_name: failed"
$ws.Range("D2").Value = "33.33/100.0"

# Row 3 (student 308418367): score moved from 0.0/100.0 to a perfect
# 100.0/100.0, and the review notes are cleared (no failures).
$ws.Range("B3").Value = "100.0 / 100.0"
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = "100.0/100.0"
